# Applies the "Add AutoCAD TakeOff Wizard / Unit 102 data" edit:
#   - Walls sheet: clears a handful of empty placeholder cells, fills in the
#     blank E5 orientation, appends the new Unit 102 (Z-102) wall rows, and
#     slides the footnote row down from row 8 to row 14.
#   - Openings sheet: clears an empty placeholder cell, appends the new
#     Unit 102 window rows, and slides the footnote row down from row 6 to
#     row 9.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Walls sheet
# ---------------------------------------------------------------------
$walls = $wb.Worksheets.Item("Walls")

# Drop the stray empty inline-string cells left over in column H / E.
$walls.Range("H2").ClearContents()
$walls.Range("H3").ClearContents()
$walls.Range("H4").ClearContents()
$walls.Range("E6").ClearContents()
$walls.Range("H6").ClearContents()
$walls.Range("E7").ClearContents()
$walls.Range("H7").ClearContents()

# Row 5 was missing its orientation - it's a West-facing demising wall.
$walls.Range("E5").Value = "West"

# The footnote currently sitting in row 8 is about to be pushed down to
# row 14 (after the new Unit 102 rows are appended). Copy its (italic,
# grey) formatting now, while it is still intact, and paste it onto the
# new footnote location; row 8 itself is reused for the first new wall.
$walls.Range("A8").Copy()
$walls.Range("A14").Value = "# Gross Area = total wall area INCLUDING windows/doors. Orientation required for exterior walls."
$walls.Range("A14:H14").Merge()
$walls.Range("A14").PasteSpecial(-4122)

$walls.Range("A8:H8").UnMerge()
$walls.Range("A8").Style = "Normal"
$walls.Range("A8").Value = "W-102-N"

$walls.Range("A9").Value = "W-102-S"
$walls.Range("B9").Value = "Z-102"
$walls.Range("C9").Value = "South Wall"
$walls.Range("D9").Value = "Exterior Wall"
$walls.Range("E9").Value = "South"
$walls.Range("F9").Value = 240
$walls.Range("G9").Value = "R-21 Wood Framed Wall"

$walls.Range("A10").Value = "W-102-W"
$walls.Range("B10").Value = "Z-102"
$walls.Range("C10").Value = "West Wall"
$walls.Range("D10").Value = "Exterior Wall"
$walls.Range("E10").Value = "West"
$walls.Range("F10").Value = 240
$walls.Range("G10").Value = "R-21 Wood Framed Wall"

$walls.Range("A11").Value = "W-102-DM"
$walls.Range("B11").Value = "Z-102"
$walls.Range("C11").Value = "Demising"
$walls.Range("D11").Value = "Interior Wall"
$walls.Range("E11").Value = "East"
$walls.Range("F11").Value = 240
$walls.Range("G11").Value = "R-0 Wall"
$walls.Range("H11").Value = "Z-101"

$walls.Range("A12").Value = "W-102-RF"
$walls.Range("B12").Value = "Z-102"
$walls.Range("C12").Value = "Roof"
$walls.Range("D12").Value = "Roof"
$walls.Range("F12").Value = 900
$walls.Range("G12").Value = "R-38 Roof"

$walls.Range("A13").Value = "W-102-SL"
$walls.Range("B13").Value = "Z-102"
$walls.Range("C13").Value = "Slab"
$walls.Range("D13").Value = "Slab on Grade"
$walls.Range("F13").Value = 900
$walls.Range("G13").Value = "Slab-on-Grade"

# ---------------------------------------------------------------------
# Openings sheet
# ---------------------------------------------------------------------
$openings = $wb.Worksheets.Item("Openings")

# Drop the stray empty SHGC placeholder on the door row.
$openings.Range("G5").ClearContents()

# Same footnote shuffle as above: row 6's footnote moves to row 9, and
# row 6 becomes the first new Unit 102 opening row.
$openings.Range("A6").Copy()
$openings.Range("A9").Value = "# U-Factor and SHGC required for windows; leave blank for doors."
$openings.Range("A9:G9").Merge()
$openings.Range("A9").PasteSpecial(-4122)

$openings.Range("A6:G6").UnMerge()
$openings.Range("A6").Style = "Normal"
$openings.Range("A6").Value = "O-102-N-1"

$openings.Range("A7").Value = "O-102-S-1"
$openings.Range("B7").Value = "W-102-S"
$openings.Range("C7").Value = "South Window"
$openings.Range("D7").Value = "Window"
$openings.Range("E7").Value = 48
$openings.Range("F7").Value = 0.27
$openings.Range("G7").Value = 0.18

$openings.Range("A8").Value = "O-102-W-1"
$openings.Range("B8").Value = "W-102-W"
$openings.Range("C8").Value = "West Window"
$openings.Range("D8").Value = "Window"
$openings.Range("E8").Value = 24
$openings.Range("F8").Value = 0.27
$openings.Range("G8").Value = 0.18
